{"js": "// Switch focus from PM10 to PM2.5\n// Find the \"PM10\" mention in the Parameter Selection paragraph and replace\n// it with \"PM2.5\", leaving the rest of the sentence untouched.\nconst body = context.document.body;\n\nconst results = body.search(\"PM10\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"PM2.5\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Switch focus from PM10 to PM2.5\n# Find the \"PM10\" mention in the Parameter Selection paragraph and replace\n# it with \"PM2.5\", leaving the rest of the sentence untouched.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"PM10\"\n$find.Replacement.Text = \"PM2.5\"\n$find.Forward = $true\n$find.Wrap = 0\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute(\"PM10\", $true, $false, $false, $false, $false, $true, 1, $false, \"PM2.5\", 2)\n"}
